$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("3_3")
$ws2 = $wb.Worksheets.Item("1_8")

# --- Add the new row-5 labels first (so the new shared string for the label
#     lands right after the existing labels, matching natural authoring order) ---
$ws1.Range("A4").Copy()
$ws1.Range("A5").PasteSpecial(-4122)
$ws1.Range("A5").Value = "BT_ACL_Sniff_1dot28s_Master_0dBm"

$ws2.Range("A4").Copy()
$ws2.Range("A5").PasteSpecial(-4122)
$ws2.Range("A5").Value = "BT_ACL_Sniff_1dot28s_Master_0dBm"

# --- Sheet "3_3": numeric updates for existing rows ---
$ws1.Range("B2").Value = 0.2393
$ws1.Range("C2").Value = 0.2465
$ws1.Range("D2").Value = 0.2313
$ws1.Range("E2").Value = 0.0049

$ws1.Range("B3").Value = 8.3785
$ws1.Range("C3").Value = 8.3981
$ws1.Range("D3").Value = 8.3554
$ws1.Range("E3").Value = 0.0129

$ws1.Range("B4").Value = 0.2398
$ws1.Range("C4").Value = 0.2456
$ws1.Range("D4").Value = 0.236
$ws1.Range("E4").Value = 0.0031

# --- Sheet "3_3": new row 5 numeric values ---
$ws1.Range("B5").Value = 0.5613
$ws1.Range("C5").Value = 3.4334
$ws1.Range("D5").Value = 0.2377
$ws1.Range("E5").Value = 0.9574
$ws1.Range("F5").Value = 10

# --- Sheet "3_3": raw-sample text in column G (claims shared-string slots
#     right after the new label, in row order) ---
$ws1.Range("G2").Value = "0.234209,0.242412,0.239419,0.246225,0.239158,0.231252,0.233389,0.239798,0.246523,0.240357"
$ws1.Range("G3").Value = "8.39487,8.380977,8.361055,8.371432,8.382853,8.398062,8.383132,8.355366,8.372432,8.384548"
$ws1.Range("G4").Value = "0.23722,0.236363,0.245604,0.236003,0.243648,0.240177,0.236363,0.241015,0.240294,0.240835"
$ws1.Range("G5").Value = "0.24353,0.237699,0.249219,0.243151,0.237959,0.239996,0.244008,3.433437,0.242214,0.241971"

# --- Sheet "1_8": numeric updates for existing rows ---
$ws2.Range("B2").Value = -0.1002
$ws2.Range("C2").Value = -0.0857
$ws2.Range("D2").Value = -0.1302
$ws2.Range("E2").Value = 0.0179

$ws2.Range("B3").Value = -5.5106
$ws2.Range("C3").Value = -5.506
$ws2.Range("D3").Value = -5.5191
$ws2.Range("E3").Value = 0.004

$ws2.Range("B4").Value = -0.4344
$ws2.Range("C4").Value = -0.0807
$ws2.Range("D4").Value = -3.4452
$ws2.Range("E4").Value = 1.0038

# --- Sheet "1_8": new row 5 numeric values ---
$ws2.Range("B5").Value = -0.101
$ws2.Range("C5").Value = -0.083
$ws2.Range("D5").Value = -0.1319
$ws2.Range("E5").Value = 0.018
$ws2.Range("F5").Value = 10

# --- Sheet "1_8": raw-sample text in column G ---
$ws2.Range("G2").Value = "-0.086304,-0.08936,-0.127971,-0.094632,-0.087705,-0.130206,-0.08739,-0.123502,-0.089366,-0.08573"
$ws2.Range("G3").Value = "-5.507963,-5.519113,-5.506302,-5.511062,-5.506042,-5.510759,-5.516316,-5.510482,-5.509951,-5.507975"
$ws2.Range("G4").Value = "-0.127706,-0.089903,-0.085433,-0.129373,-0.080741,-3.445201,-0.086267,-0.123823,-0.084063,-0.091589"
$ws2.Range("G5").Value = "-0.08965,-0.125502,-0.092681,-0.082958,-0.126292,-0.085186,-0.094614,-0.131867,-0.092947,-0.088502"
